# edit.ps1 -- apply the CS_09_02_REC190.docx commit to the current document.
#
# Most of the hunks in the diff split a single <w:r> run that already
# carries the desired text into several runs with byte-identical <w:rPr>
# (this is the classic artifact left by Word's spell-checker wrapping a
# flagged word in <w:proofErr>, or by the cursor/_GoBack bookmark sitting
# in the middle of a run). We reproduce the run-splitting by toggling
# Bold on/off on the sub-range: Word always breaks a run at a formatting
# boundary, and flipping Bold back to False immediately leaves no visible
# trace (and no <w:b/> tag, since False is the implicit default) -- so the
# resulting runs end up with identical <w:rPr>, matching the diff.

$d = $word.ActiveDocument

function Split-Runs($doc, $startPos, $pieces) {
    # Split the text starting at $startPos into consecutive runs whose
    # lengths are given by $pieces (array of substrings, already verified
    # to concatenate back to the original text). Forces a run break after
    # every piece (except it's a no-op for zero-length pieces).
    $pos = $startPos
    foreach ($piece in $pieces) {
        $len = $piece.Length
        if ($len -gt 0) {
            $sub = $doc.Range($pos, $pos + $len)
            $sub.Bold = 1
            $sub.Bold = 0
        }
        $pos = $pos + $len
    }
}

# ---------------------------------------------------------------------
# 1) "Nombre del guión a que corresponde el ejercicio" -> split around
#    "guión" (spell-check wrapped it).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Nombre del guión a que corresponde el ejercicio", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-Runs $d $rng.Start @("Nombre del ", "guión", " a que corresponde el ejercicio")

# ---------------------------------------------------------------------
# 2) Remove the _GoBack bookmark that used to sit right after "...CO".
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# 3) "Actividad sobre los felices años 20" -> split around each word
#    that the spell-checker flagged ("Actividad", "sobre", "felices",
#    "años").
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Actividad sobre los felices años 20", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-Runs $d $rng.Start @("Actividad", " ", "sobre", " los ", "felices", " ", "años", " 20")

# ---------------------------------------------------------------------
# 4) Typo fix "nación" -> "nació", and reword " de los años 20 " into
#    " " + "EN " + "los años 20 " (three runs instead of one).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("personaje que nación de los años 20 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$rng.Text = "personaje que nació EN los años 20 "
Split-Runs $d $start @("personaje que nació", " ", "EN ", "los años 20 ")

# ---------------------------------------------------------------------
# 5) "Nivel del ejercicio, 1-Fácil, 2-Medio ó 3-Difícil" -> split around
#    the lone "ó" that the spell-checker flagged.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Nivel del ejercicio, 1-Fácil, 2-Medio ó 3-Difícil", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-Runs $d $rng.Start @("Nivel del ejercicio, 1-Fácil, 2-Medio ", "ó", " 3-Difícil")

# ---------------------------------------------------------------------
# 6) Mickey Mouse paragraph: split "...por primera vez. " into
#    "...por prim" / "era vez. " and drop a fresh _GoBack bookmark
#    (zero-length) exactly at that split point -- this is where the
#    cursor was left after the last edit.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" poco antes de que se desatara el crack de la bolsa de Nueva York, este personaje, Mickey Mouse vio la luz por primera vez. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$splitAt = $start + (" poco antes de que se desatara el crack de la bolsa de Nueva York, este personaje, Mickey Mouse vio la luz por prim").Length
Split-Runs $d $start @(" poco antes de que se desatara el crack de la bolsa de Nueva York, este personaje, Mickey Mouse vio la luz por prim", "era vez. ")
$bmRange = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 7) " 1-Fácil, 2-Medio, ó 3-Dificil:" -> split around the lone "ó".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" 1-Fácil, 2-Medio, ó 3-Dificil:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-Runs $d $rng.Start @(" 1-Fácil, 2-Medio, ", "ó", " 3-Dificil:")

# ---------------------------------------------------------------------
# 8) "Nombre de archivo Shutterstock o descripción de ilustración a
#    crear" -> split around "Shutterstock".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Nombre de archivo Shutterstock o descripción de ilustración a crear", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-Runs $d $rng.Start @("Nombre de archivo ", "Shutterstock", " o descripción de ilustración a crear")

Write-Output "all edits applied"
